$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = '2026-02-08 06:18:41'
$ws.Cells.Item(2, 8).NumberFormat = '@'
$ws.Cells.Item(2, 8).Value = '90%'
$ws.Cells.Item(3, 5).Value = '2026-02-08 06:18:43'
$ws.Cells.Item(3, 14).Value = '-6.9 °C 5:52 TU'
$ws.Cells.Item(3, 15).Value = '-5.5 °C'
$ws.Cells.Item(4, 5).Value = '2026-02-08 06:18:46'
$ws.Cells.Item(4, 8).NumberFormat = '@'
$ws.Cells.Item(4, 8).Value = '72%'
$ws.Cells.Item(4, 10).Value = '1001.7 hPa'
$ws.Cells.Item(4, 15).Value = '8.3 °C'
$ws.Cells.Item(5, 5).Value = '2026-02-08 06:18:48'
$ws.Cells.Item(5, 14).Value = '-5.5 °C 5:35 TU'
$ws.Cells.Item(5, 15).Value = '-4.4 °C'
$ws.Cells.Item(6, 5).Value = '2026-02-08 06:18:51'
$ws.Cells.Item(6, 8).NumberFormat = '@'
$ws.Cells.Item(6, 8).Value = '67%'
$ws.Cells.Item(6, 10).Value = '1001.4 hPa'
$ws.Cells.Item(6, 14).Value = '6.5 °C 5:55 TU'
$ws.Cells.Item(6, 15).Value = '8.4 °C'
$ws.Cells.Item(7, 5).Value = '2026-02-08 06:18:54'
$ws.Cells.Item(7, 10).Value = '1001.6 hPa'
$ws.Cells.Item(7, 14).Value = '10.3 °C 5:30 TU'
$ws.Cells.Item(7, 15).Value = '11.1 °C'
$ws.Cells.Item(8, 5).Value = '2026-02-08 06:18:56'
$ws.Cells.Item(8, 10).Value = '1001.6 hPa'
$ws.Cells.Item(8, 14).Value = '6.4 °C 5:56 TU'
$ws.Cells.Item(8, 15).Value = '8.1 °C'
$ws.Cells.Item(9, 5).Value = '2026-02-08 06:18:59'
$ws.Cells.Item(9, 8).NumberFormat = '@'
$ws.Cells.Item(9, 8).Value = '76%'
$ws.Cells.Item(10, 5).Value = '2026-02-08 06:19:01'
$ws.Cells.Item(10, 14).Value = '4.0 °C 5:59 TU'
$ws.Cells.Item(10, 15).Value = '7.1 °C'
$ws.Cells.Item(11, 5).Value = '2026-02-08 06:19:04'
$ws.Cells.Item(11, 14).Value = '0.5 °C 5:59 TU'
$ws.Cells.Item(11, 15).Value = '1.5 °C'
$ws.Cells.Item(12, 5).Value = '2026-02-08 06:19:06'
$ws.Cells.Item(12, 8).NumberFormat = '@'
$ws.Cells.Item(12, 8).Value = '75%'
$ws.Cells.Item(12, 14).Value = '7.7 °C 5:52 TU'
$ws.Cells.Item(12, 15).Value = '8.7 °C'
$ws.Cells.Item(13, 5).Value = '2026-02-08 06:19:09'
$ws.Cells.Item(13, 8).NumberFormat = '@'
$ws.Cells.Item(13, 8).Value = '92%'
$ws.Cells.Item(13, 14).Value = '-1.3 °C 5:54 TU'
$ws.Cells.Item(13, 15).Value = '0.5 °C'
$ws.Cells.Item(14, 5).Value = '2026-02-08 06:19:12'
$ws.Cells.Item(14, 15).Value = '8.2 °C'
$ws.Cells.Item(15, 5).Value = '2026-02-08 06:19:14'
$ws.Cells.Item(15, 15).Value = '6.2 °C'
$ws.Cells.Item(16, 5).Value = '2026-02-08 06:19:17'
$ws.Cells.Item(16, 8).NumberFormat = '@'
$ws.Cells.Item(16, 8).Value = '85%'
$ws.Cells.Item(17, 5).Value = '2026-02-08 06:19:20'
$ws.Cells.Item(18, 5).Value = '2026-02-08 06:19:22'
$ws.Cells.Item(18, 8).NumberFormat = '@'
$ws.Cells.Item(18, 8).Value = '78%'
$ws.Cells.Item(18, 10).Value = '1001.8 hPa'
$ws.Cells.Item(18, 14).Value = '6.2 °C 5:33 TU'
$ws.Cells.Item(18, 15).Value = '7.8 °C'
$ws.Cells.Item(19, 5).Value = '2026-02-08 06:19:25'
$ws.Cells.Item(19, 13).Value = '2.7 °C 5:39 TU'
$ws.Cells.Item(19, 15).Value = '2.3 °C'
$ws.Cells.Item(20, 5).Value = '2026-02-08 06:19:27'
$ws.Cells.Item(20, 14).Value = '-5.5 °C 5:31 TU'
$ws.Cells.Item(20, 15).Value = '-4.9 °C'
$ws.Cells.Item(21, 5).Value = '2026-02-08 06:19:30'
$ws.Cells.Item(21, 15).Value = '2.0 °C'
$ws.Cells.Item(22, 5).Value = '2026-02-08 06:19:33'
$ws.Cells.Item(22, 14).Value = '-7.8 °C 5:35 TU'
$ws.Cells.Item(22, 15).Value = '-6.6 °C'
$ws.Cells.Item(23, 5).Value = '2026-02-08 06:19:35'
$ws.Cells.Item(23, 13).Value = '-4.5 °C 5:51 TU'
$ws.Cells.Item(23, 15).Value = '-5.4 °C'
$ws.Cells.Item(24, 5).Value = '2026-02-08 06:19:38'
$ws.Cells.Item(24, 8).NumberFormat = '@'
$ws.Cells.Item(24, 8).Value = '89%'
$ws.Cells.Item(24, 14).Value = '3.8 °C 5:51 TU'
$ws.Cells.Item(24, 15).Value = '6.6 °C'
$ws.Cells.Item(25, 5).Value = '2026-02-08 06:19:40'
$ws.Cells.Item(26, 5).Value = '2026-02-08 06:19:43'
$ws.Cells.Item(26, 10).Value = '1001.2 hPa'
$ws.Cells.Item(27, 5).Value = '2026-02-08 06:19:46'
$ws.Cells.Item(28, 5).Value = '2026-02-08 06:19:48'
$ws.Cells.Item(28, 10).Value = '1001.8 hPa'
$ws.Cells.Item(28, 14).Value = '4.1 °C 5:42 TU'
$ws.Cells.Item(28, 15).Value = '5.6 °C'
$ws.Cells.Item(29, 5).Value = '2026-02-08 06:19:51'
$ws.Cells.Item(30, 5).Value = '2026-02-08 06:19:54'
$ws.Cells.Item(30, 8).NumberFormat = '@'
$ws.Cells.Item(30, 8).Value = '64%'
$ws.Cells.Item(30, 10).Value = '1001.0 hPa'
$ws.Cells.Item(30, 14).Value = '7.0 °C 5:59 TU'
$ws.Cells.Item(30, 15).Value = '9.3 °C'
$ws.Cells.Item(31, 5).Value = '2026-02-08 06:19:56'
$ws.Cells.Item(31, 8).NumberFormat = '@'
$ws.Cells.Item(31, 8).Value = '59%'
$ws.Cells.Item(31, 14).Value = '9.3 °C 5:59 TU'
$ws.Cells.Item(32, 5).Value = '2026-02-08 06:19:59'
$ws.Cells.Item(32, 8).NumberFormat = '@'
$ws.Cells.Item(32, 8).Value = '99%'
$ws.Cells.Item(32, 15).Value = '1.9 °C'
$ws.Cells.Item(33, 5).Value = '2026-02-08 06:20:01'
$ws.Cells.Item(33, 8).NumberFormat = '@'
$ws.Cells.Item(33, 8).Value = '91%'
$ws.Cells.Item(33, 10).Value = '1003.4 hPa'
$ws.Cells.Item(33, 14).Value = '-0.6 °C 5:59 TU'
$ws.Cells.Item(33, 15).Value = '0.5 °C'
$ws.Cells.Item(34, 5).Value = '2026-02-08 06:20:04'
$ws.Cells.Item(35, 5).Value = '2026-02-08 06:20:06'
$ws.Cells.Item(35, 8).NumberFormat = '@'
$ws.Cells.Item(35, 8).Value = '79%'
$ws.Cells.Item(35, 10).Value = '1002.5 hPa'
$ws.Cells.Item(35, 14).Value = '1.0 °C 5:51 TU'
$ws.Cells.Item(35, 15).Value = '3.7 °C'
$ws.Cells.Item(36, 5).Value = '2026-02-08 06:20:09'
$ws.Cells.Item(36, 8).NumberFormat = '@'
$ws.Cells.Item(36, 8).Value = '68%'
$ws.Cells.Item(36, 10).Value = '1001.6 hPa'
$ws.Cells.Item(36, 14).Value = '9.3 °C 5:58 TU'
$ws.Cells.Item(37, 5).Value = '2026-02-08 06:20:11'
$ws.Cells.Item(37, 14).Value = '2.5 °C 5:45 TU'
$ws.Cells.Item(37, 15).Value = '3.2 °C'
$ws.Cells.Item(38, 5).Value = '2026-02-08 06:20:14'
$ws.Cells.Item(38, 8).NumberFormat = '@'
$ws.Cells.Item(38, 8).Value = '83%'
$ws.Cells.Item(38, 14).Value = '5.2 °C 5:43 TU'
$ws.Cells.Item(38, 15).Value = '7.5 °C'
$ws.Cells.Item(39, 5).Value = '2026-02-08 06:20:16'
$ws.Cells.Item(40, 5).Value = '2026-02-08 06:20:19'
$ws.Cells.Item(40, 15).Value = '2.6 °C'
$ws.Cells.Item(41, 5).Value = '2026-02-08 06:20:22'
$ws.Cells.Item(41, 8).NumberFormat = '@'
$ws.Cells.Item(41, 8).Value = '95%'
$ws.Cells.Item(41, 10).Value = '1001.3 hPa'
$ws.Cells.Item(41, 14).Value = '6.9 °C 5:59 TU'
$ws.Cells.Item(41, 15).Value = '8.8 °C'
$ws.Cells.Item(42, 5).Value = '2026-02-08 06:20:24'
$ws.Cells.Item(42, 8).NumberFormat = '@'
$ws.Cells.Item(42, 8).Value = '89%'
$ws.Cells.Item(42, 14).Value = '6.7 °C 5:55 TU'
$ws.Cells.Item(42, 15).Value = '9.3 °C'
$ws.Cells.Item(43, 5).Value = '2026-02-08 06:20:27'
$ws.Cells.Item(43, 14).Value = '2.0 °C 5:59 TU'
$ws.Cells.Item(43, 15).Value = '4.4 °C'
$ws.Cells.Item(44, 5).Value = '2026-02-08 06:20:30'
$ws.Cells.Item(44, 13).Value = '-5.0 °C 5:59 TU'
$ws.Cells.Item(45, 5).Value = '2026-02-08 06:20:32'
$ws.Cells.Item(45, 8).NumberFormat = '@'
$ws.Cells.Item(45, 8).Value = '67%'
$ws.Cells.Item(45, 10).Value = '1001.9 hPa'
$ws.Cells.Item(45, 14).Value = '-0.1 °C 5:56 TU'
$ws.Cells.Item(45, 15).Value = '2.7 °C'
$ws.Cells.Item(46, 5).Value = '2026-02-08 06:20:35'
$ws.Cells.Item(46, 8).NumberFormat = '@'
$ws.Cells.Item(46, 8).Value = '86%'
$ws.Cells.Item(46, 15).Value = '6.5 °C'
